$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46: Delay - Delay Time
$ws.Cells.Item(46, 1).Value = "Delay"
$ws.Cells.Item(46, 2).Value = "Delay Time"
$ws.Cells.Item(46, 3).Value = "delayTime"
$ws.Cells.Item(46, 4).Value = "Delay Time"
$ws.Cells.Item(46, 6).Value = "sec"
$ws.Cells.Item(46, 7).Value = 0.01
$ws.Cells.Item(46, 8).Value = 2
$ws.Cells.Item(46, 9).Value = 0.5
$ws.Cells.Item(46, 10).Value = 0.01
$ws.Cells.Item(46, 11).Value = "delaySec"

# Row 47: Delay - Feedback
$ws.Cells.Item(47, 1).Value = "Delay"
$ws.Cells.Item(47, 2).Value = "Feedback"
$ws.Cells.Item(47, 3).Value = "feedback"
$ws.Cells.Item(47, 4).Value = "Feedback"
$ws.Cells.Item(47, 6).Value = "%"
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 99.5
$ws.Cells.Item(47, 9).Value = 75
$ws.Cells.Item(47, 10).Value = 0.5
$ws.Cells.Item(47, 11).Value = "feedback"

# Row 48: Delay - Dry Mix Level
$ws.Cells.Item(48, 1).Value = "Delay"
$ws.Cells.Item(48, 2).Value = "Dry Mix Level"
$ws.Cells.Item(48, 3).Value = "dryLevel"
$ws.Cells.Item(48, 4).Value = "Dry Mix Level"
$ws.Cells.Item(48, 6).Value = "%"
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 100
$ws.Cells.Item(48, 9).Value = 100
$ws.Cells.Item(48, 10).Value = 1
$ws.Cells.Item(48, 11).Value = "dryLevel"

# Row 49: Delay - Wet Mix Level
$ws.Cells.Item(49, 1).Value = "Delay"
$ws.Cells.Item(49, 2).Value = "Wet Mix Level"
$ws.Cells.Item(49, 3).Value = "wetLevel"
$ws.Cells.Item(49, 4).Value = "Delay Mix Level"
$ws.Cells.Item(49, 6).Value = "%"
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 100
$ws.Cells.Item(49, 9).Value = 50
$ws.Cells.Item(49, 10).Value = 1
$ws.Cells.Item(49, 11).Value = "wetLevel"

[void]$ws.Range("A54").Select()
